# Update "想去人数" (column F) values across the four worksheets to match
# the refreshed data snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$vals1 = @{
    4  = 513
    5  = 2398
    8  = 84
    9  = 1690
    10 = 1690
    11 = 1395
    15 = 30
    16 = 846
    17 = 116
    18 = 188
    20 = 7496
    21 = 8475
    22 = 57
    24 = 419
    30 = 11
    33 = 1500
    34 = 7
    35 = 258
    36 = 243
    37 = 26
    40 = 790
    43 = 368
    45 = 220
    47 = 208
    49 = 30
}
foreach ($row in $vals1.Keys) {
    $ws1.Range("F$row").Value = $vals1[$row]
}

# 演出 (Sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$vals2 = @{
    19 = 313
}
foreach ($row in $vals2.Keys) {
    $ws2.Range("F$row").Value = $vals2[$row]
}

# 本地生活 (Sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$vals3 = @{
    3 = 2654
    4 = 301
}
foreach ($row in $vals3.Keys) {
    $ws3.Range("F$row").Value = $vals3[$row]
}

# 全部类型 (Sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$vals4 = @{
    6  = 301
    9  = 513
    10 = 2398
    13 = 84
    14 = 1690
    15 = 1690
    16 = 1395
    18 = 30
    19 = 846
    20 = 116
    22 = 188
    25 = 7496
    26 = 8475
    27 = 57
    33 = 258
    34 = 243
    35 = 26
    39 = 790
    43 = 368
    45 = 220
    47 = 208
    49 = 313
    50 = 30
}
foreach ($row in $vals4.Keys) {
    $ws4.Range("F$row").Value = $vals4[$row]
}

$wb.Save()
